$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting on price cells whose new values would otherwise be
# auto-parsed as numbers by Excel (losing formatting such as trailing zeros).
$textCells = @("D5","D8","D9","D15","D16","D18","D20","D22","D23","D25","D26","D27","D29","D36","D37","D38","D40","D42","D43","D44","D45","D47","D48","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values row by row.
# Row 2
$ws.Range("D2").Value = "27.482.93"
$ws.Range("E2").Value = "  -0.94%  "

# Row 3
$ws.Range("D3").Value = "1.617.91"
$ws.Range("E3").Value = "  -1.75%  "

# Row 4
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
$ws.Range("D5").Value = "211.21"
$ws.Range("E5").Value = "  -1.12%  "

# Row 6
$ws.Range("E6").Value = "  -1.44%  "

# Row 7
$ws.Range("E7").Value = "  +0.02%  "

# Row 8
$ws.Range("D8").Value = "22.87"
$ws.Range("E8").Value = "  -1.79%  "

# Row 9
$ws.Range("D9").Value = "0.261"
$ws.Range("E9").Value = "  +0.25%  "

# Row 10
$ws.Range("E10").Value = "  -0.30%  "

# Row 11
$ws.Range("E11").Value = "  -0.54%  "

# Row 12
$ws.Range("D12").Value = "1.846.24"
$ws.Range("E12").Value = "  -1.70%  "

# Row 13
$ws.Range("D13").Value = "1.619.66"
$ws.Range("E13").Value = "  -1.55%  "

# Row 14
$ws.Range("E14").Value = "  -0.22%  "

# Row 15
$ws.Range("D15").Value = "0.550"
$ws.Range("E15").Value = "  -2.63%  "

# Row 16
$ws.Range("D16").Value = "64.95"
$ws.Range("E16").Value = "  +0.99%  "

# Row 17
$ws.Range("D17").Value = "27.456.79"
$ws.Range("E17").Value = "  -0.99%  "

# Row 18
$ws.Range("D18").Value = "233.23"
$ws.Range("E18").Value = "  +0.46%  "

# Row 19
$ws.Range("D19").Value = "0.0₃0719"
$ws.Range("E19").Value = "  -0.96%  "

# Row 20
$ws.Range("D20").Value = "7.55"
$ws.Range("E20").Value = "  -1.88%  "

# Row 21
$ws.Range("E21").Value = "  +0.11%  "

# Row 22
$ws.Range("D22").Value = "4.30"
$ws.Range("E22").Value = "  -0.72%  "

# Row 23
$ws.Range("D23").Value = "10.16"
$ws.Range("E23").Value = "  +0.41%  "

# Row 24
$ws.Range("E24").Value = "  +6.04%  "

# Row 25
$ws.Range("D25").Value = "150.76"
$ws.Range("E25").Value = "  +0.49%  "

# Row 26
$ws.Range("D26").Value = "6.86"
$ws.Range("E26").Value = "  -1.72%  "

# Row 27
$ws.Range("D27").Value = "0.111"
$ws.Range("E27").Value = "  -1.01%  "

# Row 28
$ws.Range("E28").Value = "  +0.07%  "

# Row 29
$ws.Range("D29").Value = "15.56"
$ws.Range("E29").Value = "  -0.75%  "

# Row 30
$ws.Range("E30").Value = "  -1.18%  "

# Row 31
$ws.Range("E31").Value = "  -1.06%  "

# Row 32
$ws.Range("E32").Value = "  -1.30%  "

# Row 33
$ws.Range("D33").Value = "1.470.11"
$ws.Range("E33").Value = "  +1.55%  "

# Row 35
$ws.Range("E35").Value = "  -3.41%  "

# Row 36
$ws.Range("D36").Value = "2.33"
$ws.Range("E36").Value = "  -0.42%  "

# Row 37
$ws.Range("D37").Value = "0.954"
$ws.Range("E37").Value = "  +6.25%  "

# Row 38
$ws.Range("D38").Value = "0.558"
$ws.Range("E38").Value = "  -2.40%  "

# Row 39
$ws.Range("E39").Value = "  -0.61%  "

# Row 40
$ws.Range("D40").Value = "0.860"
$ws.Range("E40").Value = "  -2.94%  "

# Row 42
$ws.Range("D42").Value = "68.00"
$ws.Range("E42").Value = "  +2.41%  "

# Row 43
$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").Value = "0.984"
$ws.Range("E43").Value = "  -4.63%  "

# Row 44
$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D44").Value = "2.20"
$ws.Range("E44").Value = "  -2.24%  "

# Row 45
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "5.27"
$ws.Range("E45").Value = "  -7.62%  "

# Row 46
$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value = "1.757.38"
$ws.Range("E46").Value = "  -1.67%  "

# Row 47
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").Value = "1.73"
$ws.Range("E47").Value = "  +0.70%  "

# Row 48
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "86.59"
$ws.Range("E48").Value = "  +0.03%  "

# Row 49
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₆0105"
$ws.Range("E49").Value = "  -2.62%  "

# Row 50
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "0.101"
$ws.Range("E50").Value = "  +1.46%  "

# Row 51
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "7.69"
$ws.Range("E51").Value = "  -1.23%  "
